# Fixed naive component forecaster bug - Presentation state 11.02.
# Update the quarter-over-quarter naive error values on Sheet1 with the
# corrected figures (columns B:K, rows 24-52).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(24, 11).Value = -17.20414195516174
$ws.Cells.Item(25, 10).Value = -16.98434150865444
$ws.Cells.Item(25, 11).Value = -3.386369312788844
$ws.Cells.Item(26, 9).Value = -17.22295426235204
$ws.Cells.Item(26, 10).Value = -3.624982066486445
$ws.Cells.Item(26, 11).Value = -0.9294641821589158
$ws.Cells.Item(27, 8).Value = -17.25689497300509
$ws.Cells.Item(27, 9).Value = -3.658922777139496
$ws.Cells.Item(27, 10).Value = -0.9634048928119666
$ws.Cells.Item(27, 11).Value = 0.3458441452769406
$ws.Cells.Item(28, 7).Value = -17.13529955079305
$ws.Cells.Item(28, 8).Value = -3.537327354927458
$ws.Cells.Item(28, 9).Value = -0.8418094705999279
$ws.Cells.Item(28, 10).Value = 0.4674395674889793
$ws.Cells.Item(28, 11).Value = 1.988799384936469
$ws.Cells.Item(29, 6).Value = -17.39412186950338
$ws.Cells.Item(29, 7).Value = -3.79614967363779
$ws.Cells.Item(29, 8).Value = -1.10063178931026
$ws.Cells.Item(29, 9).Value = 0.2086172487786471
$ws.Cells.Item(29, 10).Value = 1.729977066226137
$ws.Cells.Item(29, 11).Value = -4.035764575904229
$ws.Cells.Item(30, 5).Value = -17.46219349820464
$ws.Cells.Item(30, 6).Value = -3.864221302339042
$ws.Cells.Item(30, 7).Value = -1.168703418011513
$ws.Cells.Item(30, 8).Value = 0.140545620077394
$ws.Cells.Item(30, 9).Value = 1.661905437524884
$ws.Cells.Item(30, 10).Value = -4.103836204605481
$ws.Cells.Item(30, 11).Value = 2.803854532616474
$ws.Cells.Item(31, 4).Value = -17.61617421205837
$ws.Cells.Item(31, 5).Value = -4.018202016192776
$ws.Cells.Item(31, 6).Value = -1.322684131865246
$ws.Cells.Item(31, 7).Value = -0.01343509377633867
$ws.Cells.Item(31, 8).Value = 1.507924723671151
$ws.Cells.Item(31, 9).Value = -4.257816918459215
$ws.Cells.Item(31, 10).Value = 2.649873818762741
$ws.Cells.Item(31, 11).Value = 0.2395698154090965
$ws.Cells.Item(32, 3).Value = -18.09929231679948
$ws.Cells.Item(32, 4).Value = -4.501320120933885
$ws.Cells.Item(32, 5).Value = -1.805802236606356
$ws.Cells.Item(32, 6).Value = -0.4965531985174484
$ws.Cells.Item(32, 7).Value = 1.024806618930041
$ws.Cells.Item(32, 8).Value = -4.740935023200324
$ws.Cells.Item(32, 9).Value = 2.166755714021631
$ws.Cells.Item(32, 10).Value = -0.2435482893320133
$ws.Cells.Item(32, 11).Value = -1.938804334703323
$ws.Cells.Item(33, 2).Value = -20.4108400473813
$ws.Cells.Item(33, 3).Value = -6.812867851515707
$ws.Cells.Item(33, 4).Value = -4.117349967188177
$ws.Cells.Item(33, 5).Value = -2.80810092909927
$ws.Cells.Item(33, 6).Value = -1.286741111651781
$ws.Cells.Item(33, 7).Value = -7.052482753782146
$ws.Cells.Item(33, 8).Value = -0.1447920165601908
$ws.Cells.Item(33, 9).Value = -2.555096019913835
$ws.Cells.Item(33, 10).Value = -4.250352065285145
$ws.Cells.Item(33, 11).Value = -1.311034623099504
$ws.Cells.Item(34, 2).Value = -4.980277842704087
$ws.Cells.Item(34, 3).Value = -2.284759958376557
$ws.Cells.Item(34, 4).Value = -0.9755109202876501
$ws.Cells.Item(34, 5).Value = 0.5458488971598395
$ws.Cells.Item(34, 6).Value = -5.219892744970526
$ws.Cells.Item(34, 7).Value = 1.68779799225143
$ws.Cells.Item(34, 8).Value = -0.7225060111022149
$ws.Cells.Item(34, 9).Value = -2.417762056473524
$ws.Cells.Item(34, 10).Value = 0.5215553857121161
$ws.Cells.Item(34, 11).Value = -0.8573456354163971
$ws.Cells.Item(35, 2).Value = -1.909452872482039
$ws.Cells.Item(35, 3).Value = -0.6002038343931317
$ws.Cells.Item(35, 4).Value = 0.9211559830543579
$ws.Cells.Item(35, 5).Value = -4.844585659076007
$ws.Cells.Item(35, 6).Value = 2.063105078145948
$ws.Cells.Item(35, 7).Value = -0.3471989252076966
$ws.Cells.Item(35, 8).Value = -2.042454970579006
$ws.Cells.Item(35, 9).Value = 0.8968624716066345
$ws.Cells.Item(35, 10).Value = -0.4820385495218787
$ws.Cells.Item(35, 11).Value = 1.022680634228276
$ws.Cells.Item(36, 2).Value = -0.3749684946957029
$ws.Cells.Item(36, 3).Value = 1.146391322751787
$ws.Cells.Item(36, 4).Value = -4.619350319378579
$ws.Cells.Item(36, 5).Value = 2.288340417843377
$ws.Cells.Item(36, 6).Value = -0.1219635855102677
$ws.Cells.Item(36, 7).Value = -1.817219630881577
$ws.Cells.Item(36, 8).Value = 1.122097811304063
$ws.Cells.Item(36, 9).Value = -0.2568032098244498
$ws.Cells.Item(36, 10).Value = 1.247915973925705
$ws.Cells.Item(36, 11).Value = 0.9211944755864938
$ws.Cells.Item(37, 2).Value = 1.14167028642729
$ws.Cells.Item(37, 3).Value = -4.624071355703076
$ws.Cells.Item(37, 4).Value = 2.283619381518879
$ws.Cells.Item(37, 5).Value = -0.1266846218347649
$ws.Cells.Item(37, 6).Value = -1.821940667206074
$ws.Cells.Item(37, 7).Value = 1.117376774979566
$ws.Cells.Item(37, 8).Value = -0.261524246148947
$ws.Cells.Item(37, 9).Value = 1.243194937601208
$ws.Cells.Item(37, 10).Value = 0.9164734392619965
$ws.Cells.Item(37, 11).Value = 1.754572496573351
$ws.Cells.Item(38, 2).Value = -4.896022371537698
$ws.Cells.Item(38, 3).Value = 2.011668365684257
$ws.Cells.Item(38, 4).Value = -0.3986356376693871
$ws.Cells.Item(38, 5).Value = -2.093891683040697
$ws.Cells.Item(38, 6).Value = 0.8454257591449439
$ws.Cells.Item(38, 7).Value = -0.5334752619835692
$ws.Cells.Item(38, 8).Value = 0.9712439217665854
$ws.Cells.Item(38, 9).Value = 0.6445224234273743
$ws.Cells.Item(38, 10).Value = 1.482621480738728
$ws.Cells.Item(38, 11).Value = -1.197078432822523
$ws.Cells.Item(39, 2).Value = 2.832994207660627
$ws.Cells.Item(39, 3).Value = 0.4226902043069828
$ws.Cells.Item(39, 4).Value = -1.272565841064327
$ws.Cells.Item(39, 5).Value = 1.666751601121314
$ws.Cells.Item(39, 6).Value = 0.2878505799928007
$ws.Cells.Item(39, 7).Value = 1.792569763742955
$ws.Cells.Item(39, 8).Value = 1.465848265403744
$ws.Cells.Item(39, 9).Value = 2.303947322715098
$ws.Cells.Item(39, 10).Value = -0.3757525908461526
$ws.Cells.Item(39, 11).Value = 0.8596701032167943
$ws.Cells.Item(40, 2).Value = -0.03790361708925488
$ws.Cells.Item(40, 3).Value = -1.733159662460564
$ws.Cells.Item(40, 4).Value = 1.206157779725076
$ws.Cells.Item(40, 5).Value = -0.172743241403437
$ws.Cells.Item(40, 6).Value = 1.331975942346718
$ws.Cells.Item(40, 7).Value = 1.005254444007507
$ws.Cells.Item(40, 8).Value = 1.843353501318861
$ws.Cells.Item(40, 9).Value = -0.8363464122423903
$ws.Cells.Item(40, 10).Value = 0.3990762818205566
$ws.Cells.Item(40, 11).Value = 2.123380570563001
$ws.Cells.Item(41, 2).Value = -1.502432366452369
$ws.Cells.Item(41, 3).Value = 1.436885075733271
$ws.Cells.Item(41, 4).Value = 0.05798405460475808
$ws.Cells.Item(41, 5).Value = 1.562703238354913
$ws.Cells.Item(41, 6).Value = 1.235981740015702
$ws.Cells.Item(41, 7).Value = 2.074080797327056
$ws.Cells.Item(41, 8).Value = -0.6056191162341953
$ws.Cells.Item(41, 9).Value = 0.6298035778287516
$ws.Cells.Item(41, 10).Value = 2.354107866571197
$ws.Cells.Item(41, 11).Value = 2.720553059184225
$ws.Cells.Item(42, 2).Value = 2.021185630531559
$ws.Cells.Item(42, 3).Value = 0.6422846094030465
$ws.Cells.Item(42, 4).Value = 2.147003793153201
$ws.Cells.Item(42, 5).Value = 1.82028229481399
$ws.Cells.Item(42, 6).Value = 2.658381352125344
$ws.Cells.Item(42, 7).Value = -0.02131856143590682
$ws.Cells.Item(42, 8).Value = 1.21410413262704
$ws.Cells.Item(42, 9).Value = 2.938408421369485
$ws.Cells.Item(42, 10).Value = 3.304853613982513
$ws.Cells.Item(42, 11).Value = -2.170593826049543
$ws.Cells.Item(43, 2).Value = 2.59639870328499
$ws.Cells.Item(43, 3).Value = 4.101117887035145
$ws.Cells.Item(43, 4).Value = 3.774396388695934
$ws.Cells.Item(43, 5).Value = 4.612495446007288
$ws.Cells.Item(43, 6).Value = 1.932795532446037
$ws.Cells.Item(43, 7).Value = 3.168218226508984
$ws.Cells.Item(43, 8).Value = 4.892522515251429
$ws.Cells.Item(43, 9).Value = 5.258967707864457
$ws.Cells.Item(43, 10).Value = -0.2164797321675991
$ws.Cells.Item(43, 11).Value = 3.144233349489796
$ws.Cells.Item(44, 2).Value = 2.341506873006513
$ws.Cells.Item(44, 3).Value = 2.014785374667302
$ws.Cells.Item(44, 4).Value = 2.852884431978656
$ws.Cells.Item(44, 5).Value = 0.173184518417405
$ws.Cells.Item(44, 6).Value = 1.408607212480352
$ws.Cells.Item(44, 7).Value = 3.132911501222797
$ws.Cells.Item(44, 8).Value = 3.499356693835825
$ws.Cells.Item(44, 9).Value = -1.976090746196231
$ws.Cells.Item(44, 10).Value = 1.384622335461164
$ws.Cells.Item(45, 2).Value = 1.466561183616531
$ws.Cells.Item(45, 3).Value = 2.304660240927885
$ws.Cells.Item(45, 4).Value = -0.3750396726333658
$ws.Cells.Item(45, 5).Value = 0.8603830214295811
$ws.Cells.Item(45, 6).Value = 2.584687310172026
$ws.Cells.Item(45, 7).Value = 2.951132502785054
$ws.Cells.Item(45, 8).Value = -2.524314937247002
$ws.Cells.Item(45, 9).Value = 0.8363981444103927
$ws.Cells.Item(46, 2).Value = 2.012044252255831
$ws.Cells.Item(46, 3).Value = -0.6676556613054199
$ws.Cells.Item(46, 4).Value = 0.567767032757527
$ws.Cells.Item(46, 5).Value = 2.292071321499972
$ws.Cells.Item(46, 6).Value = 2.658516514113
$ws.Cells.Item(46, 7).Value = -2.816930925919056
$ws.Cells.Item(46, 8).Value = 0.5437821557383387
$ws.Cells.Item(47, 2).Value = -0.9880269291625079
$ws.Cells.Item(47, 3).Value = 0.247395764900439
$ws.Cells.Item(47, 4).Value = 1.971700053642884
$ws.Cells.Item(47, 5).Value = 2.338145246255912
$ws.Cells.Item(47, 6).Value = -3.137302193776144
$ws.Cells.Item(47, 7).Value = 0.2234108878812506
$ws.Cells.Item(48, 2).Value = 0.5533907096283328
$ws.Cells.Item(48, 3).Value = 2.277694998370778
$ws.Cells.Item(48, 4).Value = 2.644140190983806
$ws.Cells.Item(48, 5).Value = -2.83130724904825
$ws.Cells.Item(48, 6).Value = 0.5294058326091444
$ws.Cells.Item(49, 2).Value = 1.9047312492914
$ws.Cells.Item(49, 3).Value = 2.271176441904428
$ws.Cells.Item(49, 4).Value = -3.204270998127628
$ws.Cells.Item(49, 5).Value = 0.1564420835297668
$ws.Cells.Item(50, 2).Value = 2.102905296315023
$ws.Cells.Item(50, 3).Value = -3.372542143717033
$ws.Cells.Item(50, 4).Value = -0.01182906205963841
$ws.Cells.Item(51, 2).Value = -3.650852129092033
$ws.Cells.Item(51, 3).Value = -0.2901390474346385
$ws.Cells.Item(52, 2).Value = -0.2631055417942008
